# Update the "Category" column (F) for specific rows in the SAP data sheet.
# Most flagged rows move from "Operating Activities (XGBoost)" / "Financing Activities (XGBoost)"
# to "Investing Activities (XGBoost)", with one row (393) moving the opposite way
# (from "Investing Activities (XGBoost)" back to "Operating Activities (XGBoost)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Category (column F) should become "Investing Activities (XGBoost)"
$rowsToInvesting = @(
    4,5,11,13,20,35,36,37,39,42,43,46,47,49,62,68,69,74,75,76,84,86,87,88,91,93,
    99,104,107,109,110,114,125,152,154,166,167,170,175,176,178,203,208,209,211,219,223,
    224,228,234,237,244,255,258,261,265,268,274,275,276,282,285,293,295,298,300,301,303,
    304,308,309,312,317,320,327,333,336,348,351,355,359,363,366,369,371,375,384,385,386,
    396,401,408,412,415,420,422,423,426,429,430,432,433,435,442,447,448
)

foreach ($r in $rowsToInvesting) {
    $ws.Cells.Item($r, 6).Value = "Investing Activities (XGBoost)"
}

# Row(s) whose Category (column F) should become "Operating Activities (XGBoost)"
$rowsToOperating = @(393)

foreach ($r in $rowsToOperating) {
    $ws.Cells.Item($r, 6).Value = "Operating Activities (XGBoost)"
}
